$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.062.16"
$ws.Range("E2").Value = "  +5.99%  "

$ws.Range("D3").Value = "1.716.90"
$ws.Range("E3").Value = "  +3.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3326"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.183"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07483"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.912"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.30%  "

$ws.Range("D16").Value = "1.714.47"
$ws.Range("E16").Value = "  +3.25%  "

$ws.Range("E17").Value = "  +2.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06641"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.069"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.99%  "

$ws.Range("E23").Value = "  +3.45%  "

$ws.Range("D24").Value = "25.980.66"
$ws.Range("E24").Value = "  +5.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.463"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.487"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.307"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.77%  "

$ws.Range("D30").Value = "1.906.69"
$ws.Range("E30").Value = "  +3.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.109"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.965"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08533"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.721"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.279"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02286"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2131"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.530"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.96%  "

$ws.Range("E43").Value = "  +13.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6163"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.832"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5866"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.010"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07255"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.14%  "
